# Natmi following Dr Hou advice
# Update recalculated ligand-receptor edge statistics (Sema6d -> Kdr) in
# columns E:T for data rows 2-10 on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Columns E..T, in order, for each data row.
$cols = @("E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$rowValues = @{
    2  = @(3, 1, 30.708183, 92.124549, 0.3702499640981372, 0.3702499640981371, 3, 1, 62.351132, 187.053396, 0.9620005273240912, 0.9620005273240912, 1914.689971713156, 17232.2097454184, 0.3561806607041338, 0.3561806607041337)
    3  = @(3, 1, 30.708183, 92.124549, 0.3702499640981372, 0.3702499640981371, 3, 1, 0.8562703333333334, 2.568811, 0.01321118776478095, 0.01321118776478095, 26.294506093471, 236.650554841239, 0.004891441795603897, 0.004891441795603896)
    4  = @(3, 1, 30.708183, 92.124549, 0.3702499640981372, 0.3702499640981371, 3, 1, 1.606628666666667, 4.819886, 0.02478828491112776, 0.02478828491112776, 49.33664710904601, 444.029823981414, 0.009177861598399447, 0.009177861598399445)
    5  = @(3, 1, 12.69729066666667, 38.091872, 0.1530918131325759, 0.1530918131325759, 3, 1, 62.351132, 187.053396, 0.9620005273240912, 0.9620005273240912, 791.6904463997014, 7125.214017597312, 0.1472744049625393, 0.1472744049625393)
    6  = @(3, 1, 12.69729066666667, 38.091872, 0.1530918131325759, 0.1530918131325759, 3, 1, 0.8562703333333334, 2.568811, 0.01321118776478095, 0.01321118776478095, 10.87231331157689, 97.85081980419201, 0.002022524688545219, 0.002022524688545219)
    7  = @(3, 1, 12.69729066666667, 38.091872, 0.1530918131325759, 0.1530918131325759, 3, 1, 1.606628666666667, 4.819886, 0.02478828491112776, 0.02478828491112776, 20.39983117406578, 183.598480566592, 0.003794883481491422, 0.003794883481491422)
    8  = @(3, 1, 39.53358366666667, 118.600751, 0.476658222769287, 0.476658222769287, 3, 1, 62.351132, 187.053396, 0.9620005273240912, 0.9620005273240912, 2464.963693633377, 22184.6732427004, 0.4585454616574183, 0.4585454616574182)
    9  = @(3, 1, 39.53358366666667, 118.600751, 0.476658222769287, 0.476658222769287, 3, 1, 0.8562703333333334, 2.568811, 0.01321118776478095, 0.01321118776478095, 33.85143486411789, 304.662913777061, 0.006297221280631838, 0.006297221280631838)
    10 = @(3, 1, 39.53358366666667, 118.600751, 0.476658222769287, 0.476658222769287, 3, 1, 1.606628666666667, 4.819886, 0.02478828491112776, 0.02478828491112776, 63.51578881493177, 571.6420993343861, 0.01181553983123689, 0.01181553983123689)
}

foreach ($r in $rowValues.Keys) {
    $vals = $rowValues[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$r").Value = $vals[$i]
    }
}
